$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value, continuing the list in column A
$ws.Range("A4").Value = "block3_params.xlsx"

# Update the selection to reflect the new active cell (A5, as in the diff)
$ws.Range("A5").Select()
